$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.742.47'
$ws.Range("E2").Value = '  -2.25%  '

# Row 3
$ws.Range("D3").Value = '1.560.16'
$ws.Range("E3").Value = '  -0.38%  '

# Row 4
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.14'
$ws.Range("E5").Value = '  -1.33%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.488'
$ws.Range("E6").Value = '  -2.60%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.93'
$ws.Range("E8").Value = '  -0.03%  '

# Row 9
$ws.Range("E9").Value = '  -0.85%  '

# Row 10
$ws.Range("E10").Value = '  -1.60%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0861'
$ws.Range("E11").Value = '  -0.61%  '

# Row 12
$ws.Range("D12").Value = '1.780.41'
$ws.Range("E12").Value = '  -0.46%  '

# Row 13
$ws.Range("D13").Value = '1.561.04'
$ws.Range("E13").Value = '  -0.31%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.73'
$ws.Range("E14").Value = '  -2.41%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.511'
$ws.Range("E15").Value = '  -1.00%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '61.53'
$ws.Range("E16").Value = '  -2.98%  '

# Row 17
$ws.Range("D17").Value = '26.739.74'
$ws.Range("E17").Value = '  -2.34%  '

# Row 18
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.34'
$ws.Range("E18").Value = '  +1.19%  '

# Row 19
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '213.67'
$ws.Range("E19").Value = '  +0.39%  '

# Row 20
$ws.Range("E20").Value = '  -1.93%  '

# Row 21
$ws.Range("E21").Value = '  +0.07%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.09'
$ws.Range("E22").Value = '  -0.69%  '

# Row 23
$ws.Range("E23").Value = '  -1.89%  '

# Row 24
$ws.Range("E24").Value = '  -0.28%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.46'
$ws.Range("E25").Value = '  -0.53%  '

# Row 26
$ws.Range("E26").Value = '  +0.75%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.81'
$ws.Range("E27").Value = '  -1.05%  '

# Row 28
$ws.Range("E28").Value = '  +0.00%  '

# Row 29
$ws.Range("E29").Value = '  -0.88%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0462'
$ws.Range("E30").Value = '  -1.58%  '

# Row 31
$ws.Range("E31").Value = '  -4.25%  '

# Row 32
$ws.Range("E32").Value = '  -1.66%  '

# Row 33
$ws.Range("D33").Value = '1.383.67'
$ws.Range("E33").Value = '  +0.99%  '

# Row 34
$ws.Range("E34").Value = '  -1.63%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.54'
$ws.Range("E35").Value = '  +0.05%  '

# Row 36
$ws.Range("E36").Value = '  -1.15%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.932'
$ws.Range("E37").Value = '  -3.78%  '

# Row 38
$ws.Range("E38").Value = '  -2.72%  '

# Row 39
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.519'
$ws.Range("E39").Value = '  -2.29%  '

# Row 40
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.812'
$ws.Range("E40").Value = '  -1.26%  '

# Row 41
$ws.Range("E41").Value = '  -0.03%  '

# Row 42
$ws.Range("E42").Value = '  +2.07%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.37'
$ws.Range("E43").Value = '  +1.93%  '

# Row 44
$ws.Range("E44").Value = '  +1.29%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.77'
$ws.Range("E45").Value = '  -1.39%  '

# Row 46
$ws.Range("E46").Value = '  -1.42%  '

# Row 47
$ws.Range("D47").Value = '1.694.59'
$ws.Range("E47").Value = '  -0.36%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.36'
$ws.Range("E48").Value = '  -0.20%  '

# Row 49
$ws.Range("D49").Value = '0.0₇0981'
$ws.Range("E49").Value = '  -0.93%  '

# Row 50
$ws.Range("E50").Value = '  -0.15%  '

# Row 51
$ws.Range("E51").Value = '  -0.78%  '
